# DSCI521_rna63_final_project_pres.pptx - speaker notes update
#
# The speaker notes text attached to several slides is being shuffled
# around / rewritten:
#   - Slide 19 ("Correlation Matrix") notes are cleared out.
#   - Slide 22 ("Feature Importance") notes become the old slide-19 note
#     ("This is on the dataset after cleaning and imputing missing values.").
#   - Slide 23 ("Regression Results - Baseline") notes become the old
#     slide-22 note (the Xgboost feature-importance remark, two paragraphs).
#   - Slide 9 ("Exploratory Data Analysis"), which previously had no
#     speaker notes at all, gets the Wikipedia coefficient-of-determination
#     link that used to live on slide 23's notes.
# Slide 17's notes were (and remain) empty, so nothing to do there.

$p = $ppt.ActivePresentation

function Get-NotesBody($slide) {
    $np = $slide.NotesPage
    for ($i = 1; $i -le $np.Shapes.Count; $i++) {
        $shape = $np.Shapes.Item($i)
        if ($shape.Name -like "Notes Placeholder*") {
            return $shape
        }
    }
    return $null
}

# Slide 9: add brand-new notes (previously had none).
$slide9 = $p.Slides.Item(9)
$notes9 = Get-NotesBody $slide9
$notes9.TextFrame.TextRange.Text = "https://en.wikipedia.org/wiki/Coefficient_of_determination "

# Slide 19: clear the existing note text.
$slide19 = $p.Slides.Item(19)
$notes19 = Get-NotesBody $slide19
$notes19.TextFrame.TextRange.Text = ""

# Slide 22: replace notes with what used to be on slide 19.
$slide22 = $p.Slides.Item(22)
$notes22 = Get-NotesBody $slide22
$notes22.TextFrame.TextRange.Text = "This is on the dataset after cleaning and imputing missing values."

# Slide 23: replace notes with what used to be on slide 22.
$slide23 = $p.Slides.Item(23)
$notes23 = Get-NotesBody $slide23
$notes23.TextFrame.TextRange.Text = "`rModel has a low correlation, but is very important according to Xgboost feature importance.`rProbably due to very high number of values.`r"
